$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.078.20"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "2.946.06"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'374.35"
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("D6").Value = "'102.41"
$ws.Range("E6").Value = "  -1.86%  "

$ws.Range("D7").Value = "'0.535"
$ws.Range("E7").Value = "  -1.25%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -1.51%  "

$ws.Range("D10").Value = "'36.42"
$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("E12").Value = "  -0.25%  "

$ws.Range("D13").Value = "3.399.58"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("D14").Value = "'17.88"
$ws.Range("E14").Value = "  -2.69%  "

$ws.Range("D15").Value = "'7.35"
$ws.Range("E15").Value = "  -1.61%  "

$ws.Range("D16").Value = "2.937.95"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").Value = "'0.980"
$ws.Range("E17").Value = "  +2.06%  "

$ws.Range("D18").Value = "50.941.15"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("E19").Value = "  -5.19%  "

$ws.Range("D20").Value = "'7.19"
$ws.Range("E20").Value = "  -1.99%  "

$ws.Range("D21").Value = "'12.65"
$ws.Range("E21").Value = "  -1.52%  "

$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").Value = "'264.52"
$ws.Range("E23").Value = "  +1.36%  "

$ws.Range("D24").Value = "'68.33"
$ws.Range("E24").Value = "  -1.22%  "

$ws.Range("E25").Value = "  +2.29%  "

$ws.Range("D26").Value = "'8.53"
$ws.Range("E26").Value = "  +13.16%  "

$ws.Range("D27").Value = "'8.07"
$ws.Range("E27").Value = "  +11.16%  "

$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").Value = "'0.114"
$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.168"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").Value = "'25.65"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").Value = "'50.87"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Value = "'33.70"
$ws.Range("E34").Value = "  -2.49%  "

$ws.Range("D35").Value = "'0.0447"
$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("E36").Value = "  -3.00%  "

$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("E38").Value = "  -2.25%  "

$ws.Range("D39").Value = "'2.55"
$ws.Range("E39").Value = "  -0.54%  "

$ws.Range("E40").Value = "  +0.15%  "

$ws.Range("D41").Value = "'16.44"
$ws.Range("E41").Value = "  -4.28%  "

$ws.Range("E42").Value = "  -2.31%  "

$ws.Range("D43").Value = "'120.37"
$ws.Range("E43").Value = "  -1.56%  "

$ws.Range("D44").Value = "'0.286"
$ws.Range("E44").Value = "  +0.51%  "

$ws.Range("E45").Value = "  -4.94%  "

$ws.Range("D46").Value = "'2.04"
$ws.Range("E46").Value = "  -1.46%  "

$ws.Range("D47").Value = "'3.27"
$ws.Range("E47").Value = "  +1.96%  "

$ws.Range("D48").Value = "'2.30"
$ws.Range("E48").Value = "  -3.29%  "

$ws.Range("D49").Value = "1.974.96"
$ws.Range("E49").Value = "  -2.77%  "

$ws.Range("D50").Value = "'0.0344"
$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("D51").Value = "'5.05"
$ws.Range("E51").Value = "  -0.50%  "
